$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 559
$ws.Range("F3").Value = 257
$ws.Range("F4").Value = 586
$ws.Range("F5").Value = 1395
$ws.Range("F6").Value = 707
$ws.Range("F9").Value = 163
$ws.Range("F11").Value = 6472
$ws.Range("F12").Value = 124
$ws.Range("F13").Value = 32
$ws.Range("F14").Value = 1903
$ws.Range("F15").Value = 4807
$ws.Range("F18").Value = 5625
$ws.Range("F19").Value = 7598
$ws.Range("F21").Value = 1096
$ws.Range("F22").Value = 774
$ws.Range("F23").Value = 4084
$ws.Range("F24").Value = 584
$ws.Range("F25").Value = 32
$ws.Range("F27").Value = 237
$ws.Range("F29").Value = 1078
$ws.Range("F30").Value = 1511
$ws.Range("F31").Value = 585
$ws.Range("F32").Value = 718
$ws.Range("F33").Value = 1715
$ws.Range("F34").Value = 252
$ws.Range("F35").Value = 1960
$ws.Range("F37").Value = 52
$ws.Range("F38").Value = 1274
$ws.Range("F40").Value = 710
$ws.Range("F41").Value = 329
$ws.Range("F42").Value = 2020
$ws.Range("F43").Value = 3753
$ws.Range("F44").Value = 164
$ws.Range("F45").Value = 358
$ws.Range("F46").Value = 456
$ws.Range("F47").Value = 30
$ws.Range("F48").Value = 108
$ws.Range("F49").Value = 3975
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 1292
$ws.Range("F8").Value = 3
$ws.Range("F9").Value = 21
$ws.Range("F10").Value = 21
$ws.Range("F19").Value = 10
$ws.Range("F30").Value = 88
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 4549
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 4549
$ws.Range("F5").Value = 1292
$ws.Range("F8").Value = 257
$ws.Range("F9").Value = 586
$ws.Range("F10").Value = 3
$ws.Range("F11").Value = 1395
$ws.Range("F12").Value = 21
$ws.Range("F13").Value = 707
$ws.Range("F16").Value = 163
$ws.Range("F18").Value = 6473
$ws.Range("F20").Value = 4807
$ws.Range("F21").Value = 5625
$ws.Range("F22").Value = 5625
$ws.Range("F23").Value = 7598
$ws.Range("F24").Value = 1096
$ws.Range("F25").Value = 774
$ws.Range("F26").Value = 4084
$ws.Range("F27").Value = 584
$ws.Range("F30").Value = 1078
$ws.Range("F31").Value = 1511
$ws.Range("F32").Value = 585
$ws.Range("F33").Value = 718
$ws.Range("F34").Value = 1716
$ws.Range("F35").Value = 252
$ws.Range("F36").Value = 1960
$ws.Range("F41").Value = 710
$ws.Range("F42").Value = 329
$ws.Range("F43").Value = 88
$ws.Range("F44").Value = 3753
$ws.Range("F46").Value = 164
$ws.Range("F47").Value = 358
$ws.Range("F48").Value = 108
$ws.Range("F50").Value = 3975
